$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.371.49"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.787.56"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "33.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.74%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0689"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").Value = "2.048.08"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "1.789.67"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "34.358.19"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "169.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.96%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0525"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("D35").Value = "1.412.89"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.683"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.18%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.941"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0526"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("D48").Value = "1.949.24"
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "0.0₆0127"
$ws.Range("E51").Value = "  -2.75%  "
